$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" column (C) for rows 2-11 from date serial 45192 to 45202
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45202
}
